# REFACTOR: Bring all of Dokeza up to date including the online version.
#
# The document was previously saved out of a SharePoint document library and
# picked up three legacy "Dokeza" / SharePoint custom XML metadata parts
# (the content-type schema, the SharePoint list-form template and the
# SharePoint document-management properties stub). None of that metadata is
# used anymore now that the online version has been brought up to date, so
# this strips all non-built-in Custom XML Parts from the document.

$d = $word.ActiveDocument

# Namespaces of the legacy SharePoint / Dokeza custom XML parts that need to
# go away (customXml/item1.xml, item2.xml, item3.xml in the old package).
$legacyNamespaces = @(
    "http://schemas.microsoft.com/office/2006/metadata/contentType",
    "http://schemas.microsoft.com/sharepoint/v3/contenttype/forms",
    "http://schemas.microsoft.com/office/2006/metadata/properties"
)

$allParts = $d.CustomXMLParts
$removed = 0

if ($allParts -ne $null) {
    # Walk backwards (defensive, in case Delete re-indexes the collection)
    # and drop anything that isn't a Word built-in part (core/extended
    # document properties) and whose namespace matches the legacy schema.
    for ($i = $allParts.Count; $i -ge 1; $i--) {
        $part = $allParts.Item($i)
        if ($part -eq $null) { continue }

        $isLegacy = $legacyNamespaces -contains $part.NamespaceURI
        $isBuiltIn = $false
        try { $isBuiltIn = [bool]$part.BuiltIn } catch { $isBuiltIn = $false }

        if ($isLegacy -and -not $isBuiltIn) {
            $part.Delete()
            $removed++
        }
    }

    # Belt-and-braces: also ask directly by namespace in case indexed
    # iteration above missed anything, mirroring how this is normally
    # scripted against real Word.
    foreach ($ns in $legacyNamespaces) {
        try {
            $scoped = $allParts.SelectByNamespace($ns)
        } catch {
            $scoped = $null
        }
        if ($scoped -ne $null) {
            for ($j = $scoped.Count; $j -ge 1; $j--) {
                $p = $scoped.Item($j)
                if ($p -ne $null) {
                    $p.Delete()
                    $removed++
                }
            }
        }
    }
}

Write-Output ("Removed " + $removed + " legacy custom XML part(s); remaining CustomXMLParts.Count=" + $d.CustomXMLParts.Count)
